$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.074.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  -0.20%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'" + '1.921.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +0.35%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'" + '0.9985'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'" + '  -0.31%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'" + '321.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  -2.72%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'" + '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +0.05%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'" + '0.5054'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + '  -2.64%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'" + '0.4028'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + '  -0.72%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'" + '0.08258'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + '  -2.66%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'" + '1.110'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  -1.20%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'" + '42.02'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  -1.59%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'" + '  +2.23%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'" + '1.917.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +0.19%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'" + '6.411'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  -0.31%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'" + '7.308'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + '  -1.00%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'" + '  +0.03%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'" + '92.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  -2.76%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'" + '0.00001098'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + '  -1.31%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'" + '0.06468'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  -3.36%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'" + '18.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  +1.12%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'" + '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +0.09%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'" + '5.969'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  -0.64%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'" + '30.129.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  -0.02%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'" + '11.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + '  -0.32%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'" + '2.193'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + '  -1.95%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = "'" + '2.141.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +0.26%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'" + '22.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +4.65%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'" + '161.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + '  -0.25%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'" + '2.351'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  -2.43%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'" + '129.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  +0.36%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'" + '1.129'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +3.05%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'" + '0.1044'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + '  -2.11%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'" + '  -0.23%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'" + '3.773'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +3.78%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = "'" + '0.02455'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + '  -1.34%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = "'" + '5.412'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +4.72%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'" + '0.06454'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  -1.80%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'" + '0.2161'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  -2.18%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'" + '8.901'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + '  +1.29%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'" + '1.192'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  -2.83%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'" + '0.6416'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  -1.62%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'" + '11.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  -4.30%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'" + '1.214'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  -1.99%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'" + '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +0.07%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'" + '13.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  +0.00%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = "'" + '2.170'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +4.47%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = "'" + '0.6002'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  -2.17%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'" + '3.641'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  -2.82%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'" + '123.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  -0.63%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'" + '1.216'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  -2.16%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'" + '79.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  -0.25%  '
$ws.Range("E51").Style = "Normal"
